# Generate Report for Archive
#
# The localization-status report is regenerated: the "Ready for handoff"
# status (shared across the Overview roll-up and each locale sheet) moves
# on to "In Translation", and the Status columns narrow to fit the new
# (shorter) text, matching the CI report generator's auto-fit behaviour.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$targetColumnWidth = 12.576851254417766   # ~ AutoFit width for "In Translation"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $targetColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetColumnWidth

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2:C4").Value = $newStatus
$wsZh.Columns.Item(3).ColumnWidth = $targetColumnWidth

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2:C4").Value = $newStatus
$wsDe.Columns.Item(3).ColumnWidth = $targetColumnWidth
